$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report-generated timestamp text
$ws.Range("A24").Value = "Reporte generado a las 01:31 PM el 5/12/2018"

# Update Program Size numbers for rows 11, 13, 14, 16 (column C = Actual, column D = A la Fecha)
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 670

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 670

$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 670

$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 119
